# Edit script: refresh the stock-screener lists in Sheet1
# - Column B2 is cleared (no more "NSE:GILLANDERS")
# - Column C (support zone) gets a refreshed list of tickers, rows 2-33
# - Column E (short buildup) gets updated values for rows 3-5
# - Column A (index numbers) is extended for the new rows 27-33, matching
#   the existing numbering pattern (row - 2) and the bold/bordered style
#   used by the rest of column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear B2 (previously held "NSE:GILLANDERS")
$ws.Range("B2").Value = ""

# Updated "support Zone" column (C), rows 2-26 replaced, 27-33 newly added
$cValues = @{
    2  = "NSE:APTUS"
    3  = "NSE:CENTURYPLY"
    4  = "NSE:CLEAN"
    5  = "NSE:CLSEL"
    6  = "NSE:DEVYANI"
    7  = "NSE:DWARKESH"
    8  = "NSE:ENDURANCE"
    9  = "NSE:EXIDEIND"
    10 = "NSE:GALAXYSURF"
    11 = "NSE:GOACARBON"
    12 = "NSE:HDFCMOMENT"
    13 = "NSE:HEMIPROP"
    14 = "NSE:JUBLFOOD"
    15 = "NSE:JYOTHYLAB"
    16 = "NSE:KBCGLOBAL"
    17 = "NSE:KOPRAN"
    18 = "NSE:KTKBANK"
    19 = "NSE:LUMAXIND"
    20 = "NSE:LUXIND"
    21 = "NSE:MAXESTATES"
    22 = "NSE:MRPL"
    23 = "NSE:NITINSPIN"
    24 = "NSE:NYKAA"
    25 = "NSE:OIL"
    26 = "NSE:ONMOBILE"
    27 = "NSE:ORIENTCER"
    28 = "NSE:PAGEIND"
    29 = "NSE:PANAMAPET"
    30 = "NSE:POWERGRID"
    31 = "NSE:REPRO"
    32 = "NSE:RML"
    33 = "NSE:SALASAR"
}

foreach ($row in $cValues.Keys) {
    $ws.Range("C$row").Value = $cValues[$row]
}

# Updated "Short buildup" column (E) values
$ws.Range("E3").Value = "NSE:GAIL"
$ws.Range("E4").Value = "NSE:GUJGASLTD"
$ws.Range("E5").Value = "NSE:MGL"

# Extend column A (index) and copy formatting for the new rows 27-33
for ($row = 27; $row -le 33; $row++) {
    $ws.Range("A$row").Value = $row - 2
    $ws.Range("A26").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
}

$ws.Range("A1").Select()
